$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

for ($row = 2; $row -le 19; $row++) {
    $display = $ws.Cells.Item($row, 3).Value2
    $ws.Cells.Item($row, 4).Value = $display
}
